# Apply scheduled runner updates to Typhon Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 3000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 3000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 3000
$ws.Range("N13").Value = -3338

$ws.Range("H43").Value = 319.5
$ws.Range("I43").Value = 250.25
$ws.Range("J43").Value = 365.66666
$ws.Range("K43").Value = 250.25
$ws.Range("L43").Value = 365.66666
$ws.Range("M43").Value = -181.25
$ws.Range("N43").Value = -503.66666

$ws.Range("H98").Value = 1002.8889
$ws.Range("I98").Value = 1035
$ws.Range("J98").Value = 746
$ws.Range("K98").Value = 1035
$ws.Range("L98").Value = 746
$ws.Range("M98").Value = 463
$ws.Range("N98").Value = -3742

$ws.Range("H113").Value = 37040972
$ws.Range("I113").Value = 76926456
$ws.Range("J113").Value = 4453.4287
$ws.Range("K113").Value = 76926456
$ws.Range("L113").Value = 4453.4287
$ws.Range("M113").Value = -76923202
$ws.Range("N113").Value = -10961.4287

$ws.Range("H122").Value = 1002.8889
$ws.Range("I122").Value = 1035
$ws.Range("J122").Value = 746
$ws.Range("K122").Value = 3105
$ws.Range("L122").Value = 2238
$ws.Range("M122").Value = -655
$ws.Range("N122").Value = -7138

$ws.Range("H136").Value = 50000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 50000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 50000
$ws.Range("N136").Value = -60200

$ws.Range("H137").Value = 1313.1945
$ws.Range("I137").Value = 1202.5
$ws.Range("J137").Value = 1866.6666
$ws.Range("K137").Value = 3607.5
$ws.Range("L137").Value = 5599.9998
$ws.Range("M137").Value = -1057.5
$ws.Range("N137").Value = -10699.9998

$ws.Range("H138").Value = 1617.2697
$ws.Range("I138").Value = 609.425
$ws.Range("J138").Value = 2440
$ws.Range("K138").Value = 1828.275
$ws.Range("L138").Value = 7320
$ws.Range("M138").Value = 3311.725
$ws.Range("N138").Value = -17600

$ws.Range("H141").Value = 1341
$ws.Range("I141").Value = 1230.2084
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 3690.6252
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = 1489.3748

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H74").Value = 3064.8262
$ws.Range("I74").Value = 3250.125
$ws.Range("J74").Value = 2641.2856
$ws.Range("K74").Value = 3250.125
$ws.Range("L74").Value = 2641.2856
$ws.Range("M74").Value = -2376.125

$ws.Range("H77").Value = 3064.8262
$ws.Range("I77").Value = 3250.125
$ws.Range("J77").Value = 2641.2856
$ws.Range("K77").Value = 16250.625
$ws.Range("L77").Value = 13206.428
$ws.Range("M77").Value = -11882.625

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H110").Value = 798
$ws.Range("I110").Value = 857.6
$ws.Range("J110").Value = 500
$ws.Range("K110").Value = 857.6
$ws.Range("L110").Value = 500
$ws.Range("M110").Value = 1187.4

$ws.Range("H132").Value = 16677.324
$ws.Range("I132").Value = 2104.25
$ws.Range("J132").Value = 51652.7
$ws.Range("K132").Value = 6312.75
$ws.Range("L132").Value = 154958.1
$ws.Range("M132").Value = -3782.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1553
$ws.Range("I99").Value = 1270.0555
$ws.Range("J99").Value = 2401.8333
$ws.Range("K99").Value = 1270.0555
$ws.Range("L99").Value = 2401.8333
$ws.Range("M99").Value = 227.9445000000001

$ws.Range("H100").Value = 29428.666
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 29428.666
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 29428.666
$ws.Range("N100").Value = -31592.666

$ws.Range("H107").Value = 1326
$ws.Range("I107").Value = 1122.3077
$ws.Range("J107").Value = 1855.6
$ws.Range("K107").Value = 1122.3077
$ws.Range("L107").Value = 1855.6
$ws.Range("M107").Value = 797.6922999999999
$ws.Range("N107").Value = -5695.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 7
$ws.Range("I10").Value = 7
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 7
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 132

$ws.Range("H31").Value = 11825.059
$ws.Range("I31").Value = 17029.7
$ws.Range("J31").Value = 4389.857
$ws.Range("K31").Value = 17029.7
$ws.Range("L31").Value = 4389.857
$ws.Range("M31").Value = -16734.7
$ws.Range("N31").Value = -4979.857

$ws.Range("H34").Value = 11825.059
$ws.Range("I34").Value = 17029.7
$ws.Range("J34").Value = 4389.857
$ws.Range("K34").Value = 17029.7
$ws.Range("L34").Value = 4389.857
$ws.Range("M34").Value = -16827.7
$ws.Range("N34").Value = -4793.857

$ws.Range("H111").Value = 36660.4
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 36660.4
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 36660.4
$ws.Range("N111").Value = -44840.4

$ws.Range("H132").Value = 13976.341
$ws.Range("I132").Value = 16723.705
$ws.Range("J132").Value = 4635.3
$ws.Range("K132").Value = 50171.11500000001
$ws.Range("L132").Value = 13905.9
$ws.Range("M132").Value = -47641.11500000001

$ws.Range("H134").Value = 941.2857
$ws.Range("I134").Value = 856.3889
$ws.Range("J134").Value = 1094.1
$ws.Range("K134").Value = 2569.1667
$ws.Range("L134").Value = 3282.3
$ws.Range("M134").Value = -34.16670000000022
$ws.Range("N134").Value = -8352.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 243.5
$ws.Range("I8").Value = 243.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 730.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -591.5

$ws.Range("H39").Value = 4000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 4000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 12000
$ws.Range("N39").Value = -12588

$ws.Range("H129").Value = 228098.05
$ws.Range("I129").Value = 642.5
$ws.Range("J129").Value = 358072.66
$ws.Range("K129").Value = 1927.5
$ws.Range("L129").Value = 1074217.98
$ws.Range("M129").Value = 3072.5
$ws.Range("N129").Value = -1084217.98

$ws.Range("H131").Value = 806.4400000000001
$ws.Range("I131").Value = 730
$ws.Range("J131").Value = 807.2121
$ws.Range("K131").Value = 2190
$ws.Range("L131").Value = 2421.6363
$ws.Range("M131").Value = 2850
$ws.Range("N131").Value = -12501.6363

$ws.Range("H132").Value = 1984.5385
$ws.Range("I132").Value = 1700
$ws.Range("J132").Value = 2008.25
$ws.Range("K132").Value = 15300
$ws.Range("L132").Value = 18074.25
$ws.Range("M132").Value = -12770
$ws.Range("N132").Value = -23134.25

$ws.Range("H137").Value = 22224772
$ws.Range("I137").Value = 1196.5
$ws.Range("J137").Value = 37040490
$ws.Range("K137").Value = 3589.5
$ws.Range("L137").Value = 111121470
$ws.Range("M137").Value = 1510.5
$ws.Range("N137").Value = -111131670

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 550000
$ws.Range("I3").Value = 1000000
$ws.Range("J3").Value = 100000
$ws.Range("K3").Value = 1000000
$ws.Range("L3").Value = 100000
$ws.Range("M3").Value = -999884
$ws.Range("N3").Value = -100232

$ws.Range("H7").Value = 3738358
$ws.Range("I7").Value = 5000000
$ws.Range("J7").Value = 2602880.2
$ws.Range("K7").Value = 5000000
$ws.Range("L7").Value = 2602880.2
$ws.Range("M7").Value = -4999888
$ws.Range("N7").Value = -2603104.2

$ws.Range("H8").Value = 3738358
$ws.Range("I8").Value = 5000000
$ws.Range("J8").Value = 2602880.2
$ws.Range("K8").Value = 5000000
$ws.Range("L8").Value = 2602880.2
$ws.Range("M8").Value = -4999861
$ws.Range("N8").Value = -2603158.2

$ws.Range("H113").Value = 2894.037
$ws.Range("I113").Value = 2658.75
$ws.Range("J113").Value = 3236.2727
$ws.Range("K113").Value = 2658.75
$ws.Range("L113").Value = 3236.2727
$ws.Range("M113").Value = -488.75
$ws.Range("N113").Value = -7576.2727

$ws.Range("H122").Value = 66668816
$ws.Range("I122").Value = 37038344
$ws.Range("J122").Value = 90911930
$ws.Range("K122").Value = 111115032
$ws.Range("L122").Value = 272735790
$ws.Range("M122").Value = -111112582
$ws.Range("N122").Value = -272740690

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 5000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 5000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 5000
$ws.Range("N25").Value = -5460

$ws.Range("H61").Value = 5767.6
$ws.Range("I61").Value = 3464.7273
$ws.Range("J61").Value = 8582.223
$ws.Range("K61").Value = 3464.7273
$ws.Range("L61").Value = 8582.223
$ws.Range("M61").Value = -3262.7273
$ws.Range("N61").Value = -8986.223

$ws.Range("H113").Value = 5767.6
$ws.Range("I113").Value = 3464.7273
$ws.Range("J113").Value = 8582.223
$ws.Range("K113").Value = 3464.7273
$ws.Range("L113").Value = 8582.223
$ws.Range("M113").Value = -1294.7273
$ws.Range("N113").Value = -12922.223

$ws.Range("H122").Value = 936530.3
$ws.Range("I122").Value = 2804061.8
$ws.Range("J122").Value = 2764.5715
$ws.Range("K122").Value = 8412185.399999999
$ws.Range("L122").Value = 8293.7145
$ws.Range("M122").Value = -8409735.399999999
$ws.Range("N122").Value = -13193.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 1750
$ws.Range("I30").Value = 1750
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 1750
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -1643

$ws.Range("H107").Value = 3497271.2
$ws.Range("I107").Value = 914.5714
$ws.Range("J107").Value = 7576354
$ws.Range("K107").Value = 2743.7142
$ws.Range("L107").Value = 22729062
$ws.Range("M107").Value = -823.7142000000003
$ws.Range("N107").Value = -22732902

$ws.Range("H132").Value = 1557.9445
$ws.Range("I132").Value = 1087.6923
$ws.Range("J132").Value = 2780.6
$ws.Range("K132").Value = 3263.0769
$ws.Range("L132").Value = 8341.799999999999
$ws.Range("M132").Value = -733.0769
$ws.Range("N132").Value = -13401.8
